# Commit message: "changed the dataframe locator from index-based to
# combined index-label + string"
#
# Effect on the sheet:
#  - "PLK Regular price" (column I) is refreshed for the Hu7k single-tablet
#    rows (5.95 -> 6.95).
#  - For the Kratom Kulture rows, "PLK Regular price" (column I) gets new
#    values (plain prices, or price + "(xx.xx%)" discount annotations), and
#    the now-stale "PLK Percentage Tiered Prices" (column J) is cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hu7k rows: simple price bump, 5.95 -> 6.95.
$priceBumpRows = @(6, 7, 8, 10, 11, 12)
foreach ($r in $priceBumpRows) {
    $cell = $ws.Range("I" + $r)
    # NumberFormat="@" keeps this numeric-looking value stored as literal
    # text instead of Excel re-parsing it into a float; restore the
    # original (default) style right after so no visible style changes.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = "6.95"
    $cell.Style = $origStyle
}

# Kratom Kulture rows with plain updated text prices (still numeric-looking).
$plainPriceRows = @(19, 20)
$plainPriceValues = @("4.73", "4.58")
for ($idx = 0; $idx -lt $plainPriceRows.Count; $idx++) {
    $r = $plainPriceRows[$idx]
    $v = $plainPriceValues[$idx]
    $cell = $ws.Range("I" + $r)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $v
    $cell.Style = $origStyle
    $ws.Range("J" + $r).ClearContents()
}

# Kratom Kulture rows that now show price + percentage-off text.
$discountRowsLow = @(21, 22, 23, 24)
foreach ($r in $discountRowsLow) {
    $ws.Range("I" + $r).Value = "4.14 (15.25%)"
    $ws.Range("J" + $r).ClearContents()
}

$discountRowsHigh = @(25, 26, 27, 28, 29)
foreach ($r in $discountRowsHigh) {
    $ws.Range("I" + $r).Value = "18.04 (15.25%)"
    $ws.Range("J" + $r).ClearContents()
}
